$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the city value for the second data row (E3) to "PEDERNEIRAS"
$ws.Range("E3").Value = "PEDERNEIRAS"

# Scroll/select to mimic the new view: top-left visible cell D1, active selection E4
$ws.Range("E4").Select()
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
